# Commit: "Sun, May 31, 2020  7:05:03 AM"
#
# The underlying edit:
#  1) Three tables (slides 14, 15, 16) switch from the custom
#     "Table_0" table style ({8F049480-386E-4027-84EB-CFE7548736AC}) to
#     the built-in table style {EC7F865D-8177-4FE3-A6D4-D5F545EC1BD4}.
#  2) The deck's theme (ppt/theme/theme1.xml, used by the slide master)
#     is switched from the "Integral" / "Red Violet" colour scheme to
#     the standard Office colour scheme (the font scheme and format
#     scheme are already identical between the two themes, so only the
#     colour scheme actually changes visible content).

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------
$newTableStyle = "{EC7F865D-8177-4FE3-A6D4-D5F545EC1BD4}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Swap the slide master's colour scheme to the Office palette -
# Colors(1..12) map onto dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# RGB ints are OLE colour values: R + G*256 + B*65536
$officeColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
